$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (algorithm "DT") metric updates
$ws.Range("B2").Value = 94.47
$ws.Range("C2").Value = 95.9
$ws.Range("D2").Value = 93.02
$ws.Range("E2").Value = 94.63
